$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18 (ALC)
$ws.Range("H18").Value = 779.1818
$ws.Range("I18").Value = 619
$ws.Range("J18").Value = 1500
$ws.Range("K18").Value = 619
$ws.Range("L18").Value = 1500
$ws.Range("M18").Value = -335
$ws.Range("N18").Value = -2068

# Row 64 (ALC)
$ws.Range("H64").Value = 4466.625
$ws.Range("I64").Value = 4962.5
$ws.Range("J64").Value = 3474.875
$ws.Range("K64").Value = 4962.5
$ws.Range("L64").Value = 3474.875
$ws.Range("M64").Value = -4714.5
$ws.Range("N64").Value = -3970.875

# Row 67 (ALC)
$ws.Range("H67").Value = 4466.625
$ws.Range("I67").Value = 4962.5
$ws.Range("J67").Value = 3474.875
$ws.Range("K67").Value = 4962.5
$ws.Range("L67").Value = 3474.875
$ws.Range("M67").Value = -4104.5
$ws.Range("N67").Value = -5190.875

# Row 82 (ALC)
$ws.Range("H82").Value = 8200
$ws.Range("I82").Value = 3000
$ws.Range("K82").Value = 9000
$ws.Range("M82").Value = -8594

# Row 85 (ALC)
$ws.Range("H85").Value = 8200
$ws.Range("I85").Value = 3000
$ws.Range("K85").Value = 9000
$ws.Range("M85").Value = -7596

# Row 107 (ALC)
$ws.Range("H107").Value = 10417053
$ws.Range("I107").Value = 11363967
$ws.Range("K107").Value = 11363967
$ws.Range("M107").Value = -11362047

# Row 112 (ALC)
$ws.Range("H112").Value = 21164980
$ws.Range("J112").Value = 22858146
$ws.Range("L112").Value = 68574438
$ws.Range("N112").Value = -68576654

# Row 132 (ALC)
$ws.Range("H132").Value = 1956.3077
$ws.Range("I132").Value = 1582.0834
$ws.Range("J132").Value = 6447
$ws.Range("K132").Value = 4746.2502
$ws.Range("L132").Value = 19341
$ws.Range("M132").Value = -2216.2502
$ws.Range("N132").Value = -24401

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (ARM)
$ws.Range("H2").Value = 654.2954999999999
$ws.Range("I2").Value = 517.4474
$ws.Range("K2").Value = 517.4474
$ws.Range("M2").Value = -404.4474

# Row 45 (ARM)
$ws.Range("H45").Value = 6331.5
$ws.Range("I45").Value = 6490.143
$ws.Range("K45").Value = 6490.143
$ws.Range("M45").Value = -6113.143

# Row 97 (ARM)
$ws.Range("H97").Value = 1153.3478
$ws.Range("I97").Value = 1091.3
$ws.Range("K97").Value = 1091.3
$ws.Range("M97").Value = -595.3

# Row 102 (ARM)
$ws.Range("H102").Value = 5292712.5
$ws.Range("I102").Value = 7408967
$ws.Range("K102").Value = 7408967
$ws.Range("M102").Value = -7407345

# Row 116 (ARM)
$ws.Range("H116").Value = 654.2954999999999
$ws.Range("I116").Value = 517.4474
$ws.Range("K116").Value = 517.4474
$ws.Range("M116").Value = 1776.5526

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (BSM)
$ws.Range("H3").Value = 654.2954999999999
$ws.Range("I3").Value = 517.4474
$ws.Range("K3").Value = 517.4474
$ws.Range("M3").Value = -403.4474

# Row 19 (BSM)
$ws.Range("H19").Value = 20000
$ws.Range("J19").Value = 20000
$ws.Range("L19").Value = 20000
$ws.Range("N19").Value = -20346

# Row 94 (BSM)
$ws.Range("H94").Value = 1101.9688
$ws.Range("I94").Value = 689.7083
$ws.Range("J94").Value = 2338.75
$ws.Range("K94").Value = 689.7083
$ws.Range("L94").Value = 2338.75
$ws.Range("M94").Value = -238.7083
$ws.Range("N94").Value = -3240.75

# Row 99 (BSM)
$ws.Range("H99").Value = 90910310
$ws.Range("I99").Value = 125000824
$ws.Range("K99").Value = 125000824
$ws.Range("M99").Value = -124999326

# Row 105 (BSM)
$ws.Range("H105").Value = 9557.0625
$ws.Range("I105").Value = 13945.777
$ws.Range("K105").Value = 13945.777
$ws.Range("M105").Value = -12198.777

# Row 107 (BSM)
$ws.Range("H107").Value = 874.1875
$ws.Range("I107").Value = 1056.091
$ws.Range("J107").Value = 474
$ws.Range("K107").Value = 1056.091
$ws.Range("L107").Value = 474
$ws.Range("M107").Value = 863.9090000000001
$ws.Range("N107").Value = -4314

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 4614.2334
$ws.Range("I31").Value = 1824.8536
$ws.Range("J31").Value = 10633.421
$ws.Range("K31").Value = 1824.8536
$ws.Range("L31").Value = 10633.421
$ws.Range("M31").Value = -1529.8536
$ws.Range("N31").Value = -11223.421

# Row 34 (CRP)
$ws.Range("H34").Value = 4614.2334
$ws.Range("I34").Value = 1824.8536
$ws.Range("J34").Value = 10633.421
$ws.Range("K34").Value = 1824.8536
$ws.Range("L34").Value = 10633.421
$ws.Range("M34").Value = -1622.8536
$ws.Range("N34").Value = -11037.421

# Row 36 (CRP)
$ws.Range("H36").Value = 9000
$ws.Range("J36").Value = 9000
$ws.Range("L36").Value = 9000
$ws.Range("N36").Value = -9776

# Row 40 (CRP)
$ws.Range("H40").Value = 9000
$ws.Range("J40").Value = 9000
$ws.Range("L40").Value = 9000
$ws.Range("N40").Value = -9320

# Row 105 (CRP)
$ws.Range("H105").Value = 3228.5908
$ws.Range("I105").Value = 3258.95
$ws.Range("J105").Value = 2925
$ws.Range("K105").Value = 3258.95
$ws.Range("L105").Value = 2925
$ws.Range("M105").Value = -1511.95
$ws.Range("N105").Value = -6419

# Row 107 (CRP)
$ws.Range("H107").Value = 469.11905
$ws.Range("I107").Value = 350.30768
$ws.Range("J107").Value = 662.1875
$ws.Range("K107").Value = 350.30768
$ws.Range("L107").Value = 662.1875
$ws.Range("M107").Value = 1569.69232
$ws.Range("N107").Value = -4502.1875

$ws = $wb.Worksheets.Item("CUL")
# Row 86 (CUL)
$ws.Range("H86").Value = 893.625
$ws.Range("I86").Value = 893.625
$ws.Range("K86").Value = 2680.875
$ws.Range("M86").Value = -1494.875

# Row 89 (CUL)
$ws.Range("H89").Value = 893.625
$ws.Range("I89").Value = 893.625
$ws.Range("K89").Value = 8042.625
$ws.Range("M89").Value = -2114.625

# Row 107 (CUL)
$ws.Range("H107").Value = 336.93332
$ws.Range("I107").Value = 237.15384
$ws.Range("J107").Value = 413.2353
$ws.Range("K107").Value = 711.4615200000001
$ws.Range("L107").Value = 1239.7059
$ws.Range("M107").Value = 1208.53848
$ws.Range("N107").Value = -5079.7059

# Row 109 (CUL)
$ws.Range("H109").Value = 2100
$ws.Range("I109").Value = 400
$ws.Range("J109").Value = 2409.0908
$ws.Range("K109").Value = 1200
$ws.Range("L109").Value = 7227.2724
$ws.Range("M109").Value = -160
$ws.Range("N109").Value = -9307.2724

$ws = $wb.Worksheets.Item("GSM")
# Row 97 (GSM)
$ws.Range("H97").Value = 1542.174
$ws.Range("I97").Value = 1492.579
$ws.Range("K97").Value = 1492.579
$ws.Range("M97").Value = -996.579

# Row 113 (GSM)
$ws.Range("H113").Value = 40001424
$ws.Range("J113").Value = 1800.8462
$ws.Range("L113").Value = 1800.8462
$ws.Range("N113").Value = -6140.8462

# Row 132 (GSM)
$ws.Range("H132").Value = 2254.16
$ws.Range("I132").Value = 2116.8572
$ws.Range("K132").Value = 6350.571599999999
$ws.Range("M132").Value = -3820.571599999999

$ws = $wb.Worksheets.Item("LTW")
# Row 16 (LTW)
$ws.Range("H16").Value = 2098.5
$ws.Range("I16").Value = 2098.5
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2098.5
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1928.5
$ws.Range("N16").ClearContents()

# Row 34 (LTW)
$ws.Range("H34").Value = 7510
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 7510
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 7510
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -7854

# Row 100 (LTW)
$ws.Range("H100").Value = 1260.2
$ws.Range("I100").Value = 1144.3636
$ws.Range("J100").Value = 1578.75
$ws.Range("K100").Value = 1144.3636
$ws.Range("L100").Value = 1578.75
$ws.Range("M100").Value = -603.3635999999999
$ws.Range("N100").Value = -2660.75

$ws = $wb.Worksheets.Item("WVR")
# Row 5 (WVR)
$ws.Range("H5").Value = 5001
$ws.Range("J5").Value = 5001
$ws.Range("L5").Value = 5001
$ws.Range("N5").Value = -5225

# Row 17 (WVR)
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

# Row 18 (WVR)
$ws.Range("H18").Value = 333338340
$ws.Range("J18").Value = 333338340
$ws.Range("L18").Value = 333338340
$ws.Range("N18").Value = -333338686

# Row 34 (WVR)
$ws.Range("H34").Value = 6485
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 6485
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 6485
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -6891

# Row 37 (WVR)
$ws.Range("H37").Value = 9958.25
$ws.Range("I37").Value = 3333
$ws.Range("J37").Value = 12166.667
$ws.Range("K37").Value = 3333
$ws.Range("L37").Value = 12166.667
$ws.Range("M37").Value = -3130
$ws.Range("N37").Value = -12572.667

# Row 42 (WVR)
$ws.Range("H42").Value = 3888.889
$ws.Range("J42").Value = 9666.666999999999
$ws.Range("L42").Value = 9666.666999999999
$ws.Range("N42").Value = -10422.667

# Row 107 (WVR)
$ws.Range("H107").Value = 125000780
$ws.Range("I107").Value = 200000660
$ws.Range("J107").Value = 1012
$ws.Range("K107").Value = 600001980
$ws.Range("L107").Value = 3036
$ws.Range("M107").Value = -600000060
$ws.Range("N107").Value = -6876

Write-Output "done"